$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..L correspond to years 2010..2020
$cols = @("B","C","D","E","F","G","H","I","J","K","L")

# Row 1 (years) currently only has B1 populated (with style s=1). Replicate
# that formatting across C1:L1 first (a formats-only paste reuses the same
# shared style index instead of minting new ones), then fill in the values.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:L1").PasteSpecial(-4122) | Out-Null

$years = @(2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value2 = $years[$i]
}

# Row 2 (count) - now populated for every year; 2010-2015 are all zero
# (no observations that year), 2016-2020 carry the real counts.
$counts = @(0,0,0,0,0,0,42,42,42,41,43)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value2 = $counts[$i]
}

# Rows 3-9 (mean/std/min/25%/50%/75%/max): 2010-2015 (B-G) have zero
# observations, so those summary stats are undefined/blank. Write a
# single-quote placeholder (forces a real, present text cell) and then
# strip the formatting it implies, leaving a blank text cell with no
# explicit style - matching the "endogenous zero" blanks for those years.
$ws.Range("B3:G9").Value2 = "'"
$ws.Range("B3:G9").ClearFormats() | Out-Null

# 2016-2020 (H-L) carry the real statistics.
$stats = @{
    3 = @(5.723558153020192, 5.798152265637031, 6.117503522452244, 5.82656053680151, 5.565534543592364)
    4 = @(1.253763836537182, 1.183294736722819, 1.167959067398114, 0.9344216723656189, 0.8958563041025963)
    5 = @(3.152336191851997, 3.628306484795113, 3.322370620528491, 4.220876693802421, 4.36)
    6 = @(4.864806790115457, 4.958984573235083, 5.607525680470014, 5.155529335332706, 4.774488785043467)
    7 = @(5.674900271730575, 5.948942443930747, 6.277415476188198, 5.816603052797736, 5.468483350024568)
    8 = @(6.859348464173619, 6.900025957304678, 6.827418096368399, 6.219868229362246, 6.297854220657186)
    9 = @(8.102198711914406, 7.563904173416154, 7.895786943097561, 7.754257368628767, 7.232001944008197)
}

$hlCols = @("H","I","J","K","L")
foreach ($row in $stats.Keys) {
    $vals = $stats[$row]
    for ($i = 0; $i -lt $hlCols.Length; $i++) {
        $ws.Range($hlCols[$i] + $row).Value2 = $vals[$i]
    }
}

Write-Output "done"
